$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.Value2 = "'" + $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.658.05"
Set-TextValue $ws.Range("E2") "  -0.90%  "
Set-TextValue $ws.Range("D3") "2.497.70"
Set-TextValue $ws.Range("E3") "  -1.56%  "
Set-TextValue $ws.Range("E4") "  +0.01%  "
Set-TextValue $ws.Range("D5") "569.96"
Set-TextValue $ws.Range("E5") "  -1.53%  "
Set-TextValue $ws.Range("D6") "165.84"
Set-TextValue $ws.Range("E6") "  -2.53%  "
Set-TextValue $ws.Range("D8") "0.512"
Set-TextValue $ws.Range("E8") "  -1.22%  "
Set-TextValue $ws.Range("D9") "2.496.62"
Set-TextValue $ws.Range("E9") "  -1.54%  "
Set-TextValue $ws.Range("D10") "0.158"
Set-TextValue $ws.Range("E10") "  -2.86%  "
Set-TextValue $ws.Range("E11") "  -0.46%  "
Set-TextValue $ws.Range("D12") "0.356"
Set-TextValue $ws.Range("E12") "  +0.86%  "
Set-TextValue $ws.Range("D13") "4.93"
Set-TextValue $ws.Range("E13") "  -0.08%  "
Set-TextValue $ws.Range("D14") "2.951.09"
Set-TextValue $ws.Range("E14") "  -1.22%  "
Set-TextValue $ws.Range("D15") "69.542.56"
Set-TextValue $ws.Range("E15") "  -0.89%  "
Set-TextValue $ws.Range("E16") "  -0.87%  "
Set-TextValue $ws.Range("D17") "24.36"
Set-TextValue $ws.Range("E17") "  -3.85%  "
Set-TextValue $ws.Range("D18") "2.495.44"
Set-TextValue $ws.Range("E18") "  -1.66%  "
Set-TextValue $ws.Range("D19") "11.20"
Set-TextValue $ws.Range("E19") "  -1.51%  "
Set-TextValue $ws.Range("D20") "7.40"
Set-TextValue $ws.Range("E20") "  -6.44%  "
Set-TextValue $ws.Range("E21") "  -1.75%  "
Set-TextValue $ws.Range("E22") "  -1.60%  "
Set-TextValue $ws.Range("D23") "1.93"
Set-TextValue $ws.Range("E23") "  -4.57%  "
Set-TextValue $ws.Range("E24") "  +0.06%  "
Set-TextValue $ws.Range("D25") "70.46"
Set-TextValue $ws.Range("E25") "  +0.61%  "
Set-TextValue $ws.Range("D26") "3.90"
Set-TextValue $ws.Range("E26") "  -3.26%  "
Set-TextValue $ws.Range("B27") "Aptos"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D27") "8.67"
Set-TextValue $ws.Range("E27") "  -4.89%  "
Set-TextValue $ws.Range("B28") "WrappedeETH"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D28") "2.615.29"
Set-TextValue $ws.Range("E28") "  -1.60%  "
Set-TextValue $ws.Range("D29") "0.998"
Set-TextValue $ws.Range("E29") "  -0.59%  "
Set-TextValue $ws.Range("B30") "PEPE"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D30") "0.0₃0881"
Set-TextValue $ws.Range("E30") "  -3.91%  "
Set-TextValue $ws.Range("B31") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D31") "7.83"
Set-TextValue $ws.Range("E31") "  -1.12%  "
Set-TextValue $ws.Range("D32") "445.00"
Set-TextValue $ws.Range("E32") "  -5.29%  "
Set-TextValue $ws.Range("D33") "1.20"
Set-TextValue $ws.Range("E33") "  -6.06%  "
Set-TextValue $ws.Range("E34") "  +0.04%  "
Set-TextValue $ws.Range("E35") "  -2.99%  "
Set-TextValue $ws.Range("D36") "155.89"
Set-TextValue $ws.Range("E36") "  -1.02%  "
Set-TextValue $ws.Range("E37") "  -4.17%  "
Set-TextValue $ws.Range("D38") "19.04"
Set-TextValue $ws.Range("E38") "  +0.05%  "
Set-TextValue $ws.Range("D39") "18.21"
Set-TextValue $ws.Range("E41") "  -2.08%  "
Set-TextValue $ws.Range("E42") "  -4.14%  "
Set-TextValue $ws.Range("D43") "1.59"
Set-TextValue $ws.Range("E43") "  -1.01%  "
Set-TextValue $ws.Range("D44") "38.03"
Set-TextValue $ws.Range("E44") "  -0.74%  "
Set-TextValue $ws.Range("D45") "2.17"
Set-TextValue $ws.Range("E45") "  -7.14%  "
Set-TextValue $ws.Range("E46") "  -8.13%  "
Set-TextValue $ws.Range("D47") "139.95"
Set-TextValue $ws.Range("E47") "  -3.28%  "
Set-TextValue $ws.Range("D48") "3.45"
Set-TextValue $ws.Range("E48") "  -1.96%  "
Set-TextValue $ws.Range("D49") "0.514"
Set-TextValue $ws.Range("E49") "  -3.94%  "
Set-TextValue $ws.Range("D50") "0.0729"
Set-TextValue $ws.Range("E50") "  -1.08%  "
Set-TextValue $ws.Range("D51") "0.575"
Set-TextValue $ws.Range("E51") "  -1.62%  "
